{"js": "// Remove the \"Version ... Release Date\" header paragraph and the\n// following paragraph holding the actual version number / release date\n// (e.g. \"2025/06 ... 13 June 2025\") from the top of the README.\n// These two paragraphs sit right before the \"Exclusion of Liability\"\n// heading paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Identify the two paragraphs to remove by their text content so the\n// script is resilient to any surrounding paragraphs shifting around.\nconst toDelete = [];\nfor (const p of paragraphs.items) {\n  const text = p.text;\n  const normalized = text.replace(/\\t/g, \" \").trim();\n  const isVersionHeaderRow =\n    /^Version\\b/.test(normalized) && /Release Date$/.test(normalized);\n  const isVersionValueRow =\n    /^\\d{4}\\/\\d{2}\\b/.test(normalized) && /\\d{4}$/.test(normalized);\n  if (isVersionHeaderRow || isVersionValueRow) {\n    toDelete.push(p);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Version ... Release Date\" header paragraph and the\n# following paragraph holding the actual version number / release date\n# (e.g. \"2025/06 ... 13 June 2025\") from the top of the README.\n# These two paragraphs sit right before the \"Exclusion of Liability\"\n# heading paragraph.\n\n$d = $word.ActiveDocument\n\n$targets = @()\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    $norm = ($t -replace \"`t\", \" \").Trim()\n    if ($norm -match '^Version\\b' -and $norm -match 'Release Date$') {\n        $targets += $p\n    } elseif ($norm -match '^\\d{4}/\\d{2}\\b' -and $norm -match '\\d{4}$') {\n        $targets += $p\n    }\n}\n\n# Delete from the last match to the first so earlier paragraph\n# references are not invalidated/shifted by the deletion of later ones.\nfor ($i = $targets.Count - 1; $i -ge 0; $i--) {\n    $targets[$i].Range.Delete()\n}\n"}
